# Update "想去人数" (F column) counts across sheets 展览(1), 演出(2), 全部类型(4)
# per the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 306
$ws1.Range("F3").Value = 311
$ws1.Range("F4").Value = 79
$ws1.Range("F5").Value = 381
$ws1.Range("F6").Value = 11309
$ws1.Range("F7").Value = 650
$ws1.Range("F8").Value = 103
$ws1.Range("F18").Value = 319
$ws1.Range("F19").Value = 1250
$ws1.Range("F20").Value = 64
$ws1.Range("F21").Value = 890

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 11

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 306
$ws4.Range("F3").Value = 311
$ws4.Range("F4").Value = 79
$ws4.Range("F5").Value = 381
$ws4.Range("F6").Value = 11309
$ws4.Range("F7").Value = 650
$ws4.Range("F8").Value = 103
$ws4.Range("F18").Value = 319
$ws4.Range("F19").Value = 1250
$ws4.Range("F20").Value = 64
$ws4.Range("F21").Value = 890
$ws4.Range("F23").Value = 11
